$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7744852
$ws.Range("I32").Value = 818.3333
$ws.Range("K32").Value = 818.3333
$ws.Range("M32").Value = -492.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 71429070
$ws.Range("I125").Value = 357.5
$ws.Range("J125").Value = 166667340
$ws.Range("K125").Value = 3217.5
$ws.Range("L125").Value = 1500006060
$ws.Range("M125").Value = -757.5
$ws.Range("N125").Value = -1500010980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 726.4
$ws.Range("I129").Value = 366.16666
$ws.Range("J129").Value = 1266.75
$ws.Range("K129").Value = 1098.49998
$ws.Range("L129").Value = 3800.25
$ws.Range("M129").Value = 3901.50002
$ws.Range("N129").Value = -13800.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2082.7234
$ws.Range("I74").Value = 1129.8064
$ws.Range("J74").Value = 3929
$ws.Range("K74").Value = 1129.8064
$ws.Range("L74").Value = 3929
$ws.Range("M74").Value = -255.8063999999999
$ws.Range("N74").Value = -5677

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2082.7234
$ws.Range("I77").Value = 1129.8064
$ws.Range("J77").Value = 3929
$ws.Range("K77").Value = 5649.031999999999
$ws.Range("L77").Value = 19645
$ws.Range("M77").Value = -1281.031999999999
$ws.Range("N77").Value = -28381

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 313
$ws.Range("I25").Value = 313
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 313
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -78
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4171.4194
$ws.Range("I31").Value = 3647.3845
$ws.Range("J31").Value = 4549.8887
$ws.Range("K31").Value = 3647.3845
$ws.Range("L31").Value = 4549.8887
$ws.Range("M31").Value = -3352.3845
$ws.Range("N31").Value = -5139.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4171.4194
$ws.Range("I34").Value = 3647.3845
$ws.Range("J34").Value = 4549.8887
$ws.Range("K34").Value = 3647.3845
$ws.Range("L34").Value = 4549.8887
$ws.Range("M34").Value = -3445.3845
$ws.Range("N34").Value = -4953.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 63905.332
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 63905.332
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 63905.332
$ws.Range("N59").Value = -66195.33199999999
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 34885.195
$ws.Range("I99").Value = 73708.14
$ws.Range("K99").Value = 73708.14
$ws.Range("M99").Value = -72210.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 34885.195
$ws.Range("I126").Value = 73708.14
$ws.Range("K126").Value = 221124.42
$ws.Range("M126").Value = -218654.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 515.5714
$ws.Range("I5").Value = 361.8
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 1085.4
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = -973.4000000000001
$ws.Range("N5").Value = -2924

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 163.18182
$ws.Range("I50").Value = 136.875
$ws.Range("J50").Value = 233.33333
$ws.Range("K50").Value = 410.625
$ws.Range("L50").Value = 699.99999
$ws.Range("M50").Value = 70.375
$ws.Range("N50").Value = -1661.99999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 163.18182
$ws.Range("I53").Value = 136.875
$ws.Range("J53").Value = 233.33333
$ws.Range("K53").Value = 410.625
$ws.Range("L53").Value = 699.99999
$ws.Range("M53").Value = 70.375
$ws.Range("N53").Value = -1661.99999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2480
$ws.Range("J125").Value = 3166.6667
$ws.Range("L125").Value = 9500.000100000001
$ws.Range("N125").Value = -19340.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 515.5714
$ws.Range("I135").Value = 361.8
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 3256.2
$ws.Range("L135").Value = 8100
$ws.Range("M135").Value = -721.2000000000003
$ws.Range("N135").Value = -13170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2871.3809
$ws.Range("I126").Value = 2680.1428
$ws.Range("K126").Value = 8040.428400000001
$ws.Range("M126").Value = -5570.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3889.25
$ws.Range("I132").Value = 3945.84
$ws.Range("J132").Value = 3760.6365
$ws.Range("K132").Value = 11837.52
$ws.Range("L132").Value = 11281.9095
$ws.Range("M132").Value = -9307.52
$ws.Range("N132").Value = -16341.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2353.2144
$ws.Range("I7").Value = 2152.2727
$ws.Range("K7").Value = 2152.2727
$ws.Range("M7").Value = -2040.2727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1645.15
$ws.Range("I16").Value = 1725.0625
$ws.Range("J16").Value = 1325.5
$ws.Range("K16").Value = 1725.0625
$ws.Range("L16").Value = 1325.5
$ws.Range("M16").Value = -1555.0625
$ws.Range("N16").Value = -1665.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 548.4167
$ws.Range("I22").Value = 385.125
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 385.125
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -90.125
$ws.Range("N22").Value = -1465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 548.4167
$ws.Range("I27").Value = 385.125
$ws.Range("J27").Value = 875
$ws.Range("K27").Value = 385.125
$ws.Range("L27").Value = 875
$ws.Range("M27").Value = -278.125
$ws.Range("N27").Value = -1089

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 829.65515
$ws.Range("I46").Value = 646.25
$ws.Range("J46").Value = 899.5238000000001
$ws.Range("K46").Value = 646.25
$ws.Range("L46").Value = 899.5238000000001
$ws.Range("M46").Value = -458.25
$ws.Range("N46").Value = -1275.5238

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2430
$ws.Range("I68").Value = 1854.2858
$ws.Range("J68").Value = 3101.6667
$ws.Range("K68").Value = 1854.2858
$ws.Range("L68").Value = 3101.6667
$ws.Range("M68").Value = -1105.2858
$ws.Range("N68").Value = -4599.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2430
$ws.Range("I71").Value = 1854.2858
$ws.Range("J71").Value = 3101.6667
$ws.Range("K71").Value = 9271.429
$ws.Range("L71").Value = 15508.3335
$ws.Range("M71").Value = -5527.429
$ws.Range("N71").Value = -22996.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2191.4546
$ws.Range("I93").Value = 2011.7778
$ws.Range("K93").Value = 2011.7778
$ws.Range("M93").Value = -763.7778000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2353.2144
$ws.Range("I126").Value = 2152.2727
$ws.Range("K126").Value = 6456.8181
$ws.Range("M126").Value = -3986.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 344
$ws.Range("I107").Value = 425
$ws.Range("J107").Value = 279.2
$ws.Range("K107").Value = 1275
$ws.Range("L107").Value = 837.5999999999999
$ws.Range("M107").Value = 645
$ws.Range("N107").Value = -4677.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20318.574
$ws.Range("I132").Value = 30495.646
$ws.Range("J132").Value = 3017.55
$ws.Range("K132").Value = 91486.93799999999
$ws.Range("L132").Value = 9052.650000000001
$ws.Range("M132").Value = -88956.93799999999
$ws.Range("N132").Value = -14112.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 36038172
$ws.Range("I136").Value = 55557372
$ws.Range("J136").Value = 17546294
$ws.Range("K136").Value = 166672116
$ws.Range("L136").Value = 52638882
$ws.Range("M136").Value = -166669566
$ws.Range("N136").Value = -52643982
